$wb = $excel.ActiveWorkbook

# Map of sheet name -> column F cell updates (row -> new value)
$updates = @{
    "展览" = @{
        3  = 1245
        5  = 915
        6  = 1652
        7  = 355
        8  = 1113
        9  = 43
        10 = 93
        11 = 239
        12 = 5
        14 = 592
        15 = 110
        16 = 64
        17 = 19
        20 = 56
        21 = 622
        22 = 612
        23 = 100
        25 = 821
        27 = 31
        28 = 9
        29 = 226
        32 = 390
    }
    "演出" = @{
        5 = 4
        6 = 22
        7 = 223
        8 = 78
    }
    "本地生活" = @{
        2 = 289
    }
    "全部类型" = @{
        2  = 289
        4  = 1245
        6  = 915
        7  = 1652
        8  = 355
        9  = 1113
        10 = 43
        12 = 93
        13 = 239
        14 = 5
        16 = 592
        17 = 110
        18 = 64
        20 = 19
        23 = 4
        25 = 22
        26 = 223
        27 = 223
        28 = 56
        29 = 622
        30 = 612
        31 = 101
        33 = 821
        35 = 78
        36 = 31
        37 = 9
        38 = 226
        45 = 390
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
